$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{Row=2; A="sub-OAS30001_ses-d0129_run-01_T1w.nii.gz"; C=2}
    @{Row=3; A="sub-OAS30002_ses-d0371_T1w.nii.gz"; C=1}
    @{Row=4; A="sub-OAS30003_ses-d0558_run-01_T1w.nii.gz"; C=2}
    @{Row=5; A="sub-OAS30004_ses-d1101_run-01_T1w.nii.gz"; C=2}
    @{Row=6; A="sub-OAS30005_ses-d0143_T1w.nii.gz"; C=2}
    @{Row=7; A="sub-OAS30006_ses-d0166_T1w.nii.gz"; C=1}
    @{Row=8; A="sub-OAS30008_ses-d0061_run-01_T1w.nii.gz"; C=2}
    @{Row=9; A="sub-OAS30009_ses-d0148_run-01_T1w.nii.gz"; C=2}
    @{Row=10; A="sub-OAS30010_ses-d0068_T1w.nii.gz"; C=2}
    @{Row=11; A="sub-OAS30011_ses-d0055_T1w.nii.gz"; C=2}
    @{Row=12; A="sub-OAS30014_ses-d0196_run-01_T1w.nii.gz"; C=2}
    @{Row=13; A="sub-OAS30015_ses-d0116_run-01_T1w.nii.gz"; C=1}
    @{Row=14; A="sub-OAS30017_ses-d0054_run-01_T1w.nii.gz"; C=2}
    @{Row=15; A="sub-OAS30018_ses-d0070_run-01_T1w.nii.gz"; C=1}
    @{Row=16; A="sub-OAS30020_ses-d0092_run-02_T1w.nii.gz"; C=2}
    @{Row=17; A="sub-OAS30021_ses-d0071_T1w.nii.gz"; C=1}
    @{Row=18; A="sub-OAS30026_ses-d0048_T1w.nii.gz"; C=1}
    @{Row=19; A="sub-OAS30030_ses-d0170_T1w.nii.gz"; C=1}
    @{Row=20; A="sub-OAS30032_ses-d0262_run-02_T1w.nii.gz"; C=2}
    @{Row=21; A="sub-OAS30033_ses-d0133_run-02_T1w.nii.gz"; C=2}
    @{Row=22; A="sub-OAS30034_ses-d0044_T1w.nii.gz"; C=2}
    @{Row=23; A="sub-OAS30036_ses-d0059_T1w.nii.gz"; C=2}
    @{Row=24; A="sub-OAS30037_ses-d0154_T1w.nii.gz"; C=2}
    @{Row=25; A="sub-OAS30039_ses-d0103_T1w.nii.gz"; C=2}
    @{Row=26; A="sub-OAS30042_ses-d0067_T1w.nii.gz"; C=2}
    @{Row=27; A="sub-OAS30044_ses-d0054_T1w.nii.gz"; C=2}
    @{Row=28; A="sub-OAS30046_ses-d0072_run-01_T1w.nii.gz"; C=1}
    @{Row=29; A="sub-OAS30048_ses-d0983_run-02_T1w.nii.gz"; C=2}
    @{Row=30; A="sub-OAS30049_ses-d0013_run-01_T1w.nii.gz"; C=1}
    @{Row=31; A="sub-OAS30050_ses-d0110_T1w.nii.gz"; C=2}
    @{Row=32; A="sub-OAS30053_ses-d0428_run-01_T1w.nii.gz"; C=1}
    @{Row=33; A="sub-OAS30059_ses-d0230_run-02_T1w.nii.gz"; C=2}
    @{Row=34; A="sub-OAS30060_ses-d0074_run-01_T1w.nii.gz"; C=2}
    @{Row=35; A="sub-OAS30062_ses-d0087_run-02_T1w.nii.gz"; C=2}
    @{Row=36; A="sub-OAS30065_ses-d0548_T1w.nii.gz"; C=1}
    @{Row=37; A="sub-OAS30066_ses-d0524_T1w.nii.gz"; C=2}
    @{Row=38; A="sub-OAS30067_ses-d0057_T1w.nii.gz"; C=1}
    @{Row=39; A="sub-OAS30073_ses-d0033_run-02_T1w.nii.gz"; C=1}
    @{Row=40; A="sub-OAS30075_ses-d0143_T1w.nii.gz"; C=2}
    @{Row=41; A="sub-OAS30077_ses-d0944_T1w.nii.gz"; C=2}
    @{Row=42; A="sub-OAS30079_ses-d0019_run-01_T1w.nii.gz"; C=2}
    @{Row=43; A="sub-OAS30080_ses-d0048_T1w.nii.gz"; C=2}
    @{Row=44; A="sub-OAS30082_ses-d1700_run-02_T1w.nii.gz"; C=2}
    @{Row=45; A="sub-OAS30083_ses-d0465_run-01_T1w.nii.gz"; C=2}
    @{Row=46; A="sub-OAS30086_ses-d0000_run-01_T1w.nii.gz"; C=2}
    @{Row=47; A="sub-OAS30088_ses-d0093_run-01_T1w.nii.gz"; C=1}
    @{Row=48; A="sub-OAS30089_ses-d0001_T1w.nii.gz"; C=1}
    @{Row=49; A="sub-OAS30090_ses-d0118_T1w.nii.gz"; C=1}
    @{Row=50; A="sub-OAS30092_ses-d0636_T1w.nii.gz"; C=2}
    @{Row=51; A="sub-OAS30093_ses-d0056_T1w.nii.gz"; C=2}
)

foreach ($item in $rowData) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}
